$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.1424542995417397
$ws.Cells.Item(2, 4).Value = 0.2011023862291026
$ws.Cells.Item(2, 5).Value = 0.1778968200581872
$ws.Cells.Item(2, 6).Value = 0.6524037515926082
$ws.Cells.Item(2, 7).Value = 0.2904879055989156
$ws.Cells.Item(2, 8).Value = 0.4574657667400999
$ws.Cells.Item(2, 10).Value = 0.246557002603268
$ws.Cells.Item(2, 13).Value = 5.546992588574966
$ws.Cells.Item(2, 15).Value = 1.421277761602738

$ws.Cells.Item(3, 2).Value = 0.1329270314379301
$ws.Cells.Item(3, 4).Value = 0.2061636599317538
$ws.Cells.Item(3, 5).Value = 0.1758693561622593
$ws.Cells.Item(3, 6).Value = 0.670974825369882
$ws.Cells.Item(3, 7).Value = 0.2939821842161052
$ws.Cells.Item(3, 8).Value = 0.4652759133666251
$ws.Cells.Item(3, 10).Value = 0.2339875107096816
$ws.Cells.Item(3, 13).Value = 4.862844394437758
$ws.Cells.Item(3, 15).Value = 1.444658577017236

$ws.Cells.Item(4, 2).Value = 0.127146225585733
$ws.Cells.Item(4, 4).Value = 0.2094778388629521
$ws.Cells.Item(4, 5).Value = 0.1748423915855355
$ws.Cells.Item(4, 6).Value = 0.6833202303861121
$ws.Cells.Item(4, 7).Value = 0.2966503035784953
$ws.Cells.Item(4, 8).Value = 0.4704995684757094
$ws.Cells.Item(4, 10).Value = 0.2265272836485082
$ws.Cells.Item(4, 13).Value = 4.440771257589233
$ws.Cells.Item(4, 15).Value = 1.461006988380817

$ws.Cells.Item(5, 2).Value = 0.124807969784456
$ws.Cells.Item(5, 4).Value = 0.2108800569751637
$ws.Cells.Item(5, 5).Value = 0.1744781841268974
$ws.Cells.Item(5, 6).Value = 0.6885858796839024
$ws.Cells.Item(5, 7).Value = 0.2978673455329783
$ws.Cells.Item(5, 8).Value = 0.4727352548825365
$ws.Cells.Item(5, 10).Value = 0.2235511146765958
$ws.Cells.Item(5, 13).Value = 4.268272601079673
$ws.Cells.Item(5, 15).Value = 1.468165213375272

$ws.Cells.Item(6, 2).Value = 0.1244207630766425
$ws.Cells.Item(6, 4).Value = 0.2111160066951925
$ws.Cells.Item(6, 5).Value = 0.1744209723910046
$ws.Cells.Item(6, 6).Value = 0.6894743519688653
$ws.Cells.Item(6, 7).Value = 0.2980772232701838
$ws.Cells.Item(6, 8).Value = 0.4731129321276484
$ws.Cells.Item(6, 10).Value = 0.2230607628350327
$ws.Cells.Item(6, 13).Value = 4.239599228025355
$ws.Cells.Item(6, 15).Value = 1.469383651371857

$ws.Cells.Item(7, 2).Value = 0.1271146201467275
$ws.Cells.Item(7, 4).Value = 0.2094965408381242
$ws.Cells.Item(7, 5).Value = 0.1748372605898183
$ws.Cells.Item(7, 6).Value = 0.6833902970605763
$ws.Cells.Item(7, 7).Value = 0.296666193816371
$ws.Cells.Item(7, 8).Value = 0.4705292873484268
$ws.Cells.Item(7, 10).Value = 0.2264868881593571
$ws.Cells.Item(7, 13).Value = 4.438446901319622
$ws.Cells.Item(7, 15).Value = 1.46110152444237

$ws.Cells.Item(8, 2).Value = 0.1391550521503291
$ws.Cells.Item(8, 4).Value = 0.2028044821121142
$ws.Cells.Item(8, 5).Value = 0.1771522103618892
$ws.Cells.Item(8, 6).Value = 0.6586100511675106
$ws.Cells.Item(8, 7).Value = 0.2915832677034444
$ws.Cells.Item(8, 8).Value = 0.4600694328751516
$ws.Cells.Item(8, 10).Value = 0.2421690871316429
$ws.Cells.Item(8, 13).Value = 5.311515526824621
$ws.Cells.Item(8, 15).Value = 1.428922969451065

$ws.Cells.Item(9, 2).Value = 0.1633096474751596
$ws.Cells.Item(9, 4).Value = 0.1913324215752823
$ws.Cells.Item(9, 5).Value = 0.18344552213523
$ws.Cells.Item(9, 6).Value = 0.6175959377957909
$ws.Cells.Item(9, 7).Value = 0.2858368229487382
$ws.Cells.Item(9, 8).Value = 0.4429849308802574
$ws.Cells.Item(9, 10).Value = 0.2750064583307505
$ws.Cells.Item(9, 13).Value = 7.007702239187438
$ws.Cells.Item(9, 15).Value = 1.381850811136218

$ws.Cells.Item(10, 2).Value = 0.1813842619805683
$ws.Cells.Item(10, 4).Value = 0.1839274416348857
$ws.Cells.Item(10, 5).Value = 0.1891753285734481
$ws.Cells.Item(10, 6).Value = 0.5922202471140565
$ws.Cells.Item(10, 7).Value = 0.2842909973730059
$ws.Cells.Item(10, 8).Value = 0.4325628556206453
$ws.Cells.Item(10, 10).Value = 0.3004674241690708
$ws.Cells.Item(10, 13).Value = 8.244409026353196
$ws.Cells.Item(10, 15).Value = 1.357342429715629

$ws.Cells.Item(11, 2).Value = 0.1896777238767129
$ws.Cells.Item(11, 4).Value = 0.1807847701872802
$ws.Cells.Item(11, 5).Value = 0.192030793056837
$ws.Cells.Item(11, 6).Value = 0.5817386181956579
$ws.Cells.Item(11, 7).Value = 0.2841910502884701
$ws.Cells.Item(11, 8).Value = 0.4282928039644815
$ws.Cells.Item(11, 10).Value = 0.312355616623023
$ws.Cells.Item(11, 13).Value = 8.805039040506074
$ws.Cells.Item(11, 15).Value = 1.348446195448503

$ws.Cells.Item(12, 2).Value = 0.1928284108438874
$ws.Cells.Item(12, 4).Value = 0.1796275370309814
$ws.Cells.Item(12, 5).Value = 0.1931486157041036
$ws.Cells.Item(12, 6).Value = 0.5779246582451663
$ws.Cells.Item(12, 7).Value = 0.2842417744123509
$ws.Cells.Item(12, 8).Value = 0.4267443049057391
$ws.Cells.Item(12, 10).Value = 0.3169026404015369
$ws.Cells.Item(12, 13).Value = 9.017059504469444
$ws.Cells.Item(12, 15).Value = 1.345406766665548

$ws.Cells.Item(13, 2).Value = 0.1921494050662602
$ws.Cells.Item(13, 4).Value = 0.1798753022518724
$ws.Cells.Item(13, 5).Value = 0.1929062364919645
$ws.Cells.Item(13, 6).Value = 0.5787391198799412
$ws.Cells.Item(13, 7).Value = 0.2842268816299764
$ws.Cells.Item(13, 8).Value = 0.4270747444165295
$ws.Cells.Item(13, 10).Value = 0.3159213256708995
$ws.Cells.Item(13, 13).Value = 8.971409404023916
$ws.Cells.Item(13, 15).Value = 1.346046625712034

$ws.Cells.Item(14, 2).Value = 0.1899367303224579
$ws.Cells.Item(14, 4).Value = 0.1806889039820376
$ws.Cells.Item(14, 5).Value = 0.1921220208250318
$ws.Cells.Item(14, 6).Value = 0.5814217164591469
$ws.Cells.Item(14, 7).Value = 0.2841934377896962
$ws.Cells.Item(14, 8).Value = 0.4281640313448776
$ws.Cells.Item(14, 10).Value = 0.3127287893162816
$ws.Cells.Item(14, 13).Value = 8.822487640287306
$ws.Cells.Item(14, 15).Value = 1.348189507238573

$ws.Cells.Item(15, 2).Value = 0.1885827195819445
$ws.Cells.Item(15, 4).Value = 0.1811915441495344
$ws.Cells.Item(15, 5).Value = 0.1916464440942534
$ws.Cells.Item(15, 6).Value = 0.583085171284722
$ws.Cells.Item(15, 7).Value = 0.2841845414872353
$ws.Cells.Item(15, 8).Value = 0.4288401909844168
$ws.Cells.Item(15, 10).Value = 0.310779196000226
$ws.Cells.Item(15, 13).Value = 8.731232674352157
$ws.Cells.Item(15, 15).Value = 1.349545141284835

$ws.Cells.Item(16, 2).Value = 0.1808436885707749
$ws.Cells.Item(16, 4).Value = 0.1841373974435641
$ws.Cells.Item(16, 5).Value = 0.1889937876521586
$ws.Cells.Item(16, 6).Value = 0.5929268281725442
$ws.Cells.Item(16, 7).Value = 0.2843098239111868
$ws.Cells.Item(16, 8).Value = 0.4328514551066718
$ws.Cells.Item(16, 10).Value = 0.2996967621012487
$ws.Cells.Item(16, 13).Value = 8.207731445582567
$ws.Cells.Item(16, 15).Value = 1.357969614898366

$ws.Cells.Item(17, 2).Value = 0.1761142076767612
$ws.Cells.Item(17, 4).Value = 0.1860026822097538
$ws.Cells.Item(17, 5).Value = 0.1874307331997755
$ws.Cells.Item(17, 6).Value = 0.5992381217104636
$ws.Cells.Item(17, 7).Value = 0.284542483680255
$ws.Cells.Item(17, 8).Value = 0.4354334013119256
$ws.Cells.Item(17, 10).Value = 0.2929772486357365
$ws.Cells.Item(17, 13).Value = 7.886081661724347
$ws.Cells.Item(17, 15).Value = 1.363718600459265

$ws.Cells.Item(18, 2).Value = 0.173400646048222
$ws.Cells.Item(18, 4).Value = 0.1870967854356991
$ws.Cells.Item(18, 5).Value = 0.1865550924321937
$ws.Cells.Item(18, 6).Value = 0.6029680112172784
$ws.Cells.Item(18, 7).Value = 0.2847329751101739
$ws.Cells.Item(18, 8).Value = 0.4369627566199341
$ws.Cells.Item(18, 10).Value = 0.2891410738347275
$ws.Cells.Item(18, 13).Value = 7.700893020178853
$ws.Cells.Item(18, 15).Value = 1.367236977618887

$ws.Cells.Item(19, 2).Value = 0.1724830366925687
$ws.Cells.Item(19, 4).Value = 0.1874708677183037
$ws.Cells.Item(19, 5).Value = 0.1862626107626042
$ws.Cells.Item(19, 6).Value = 0.6042479524250552
$ws.Cells.Item(19, 7).Value = 0.2848071502373699
$ws.Cells.Item(19, 8).Value = 0.4374881538624038
$ws.Cells.Item(19, 10).Value = 0.2878471055782938
$ws.Cells.Item(19, 13).Value = 7.638159621232376
$ws.Cells.Item(19, 15).Value = 1.368464436277378

$ws.Cells.Item(20, 2).Value = 0.1766169752704059
$ws.Cells.Item(20, 4).Value = 0.1858019186496875
$ws.Cells.Item(20, 5).Value = 0.1875946974460092
$ws.Cells.Item(20, 6).Value = 0.5985559266830265
$ws.Cells.Item(20, 7).Value = 0.2845118381692089
$ws.Cells.Item(20, 8).Value = 0.4351539596772227
$ws.Cells.Item(20, 10).Value = 0.2936895720170298
$ws.Cells.Item(20, 13).Value = 7.920340875635816
$ws.Cells.Item(20, 15).Value = 1.363084660377751

$ws.Cells.Item(21, 2).Value = 0.1905863723777088
$ws.Cells.Item(21, 4).Value = 0.1804490356834805
$ws.Cells.Item(21, 5).Value = 0.1923513669066068
$ws.Cells.Item(21, 6).Value = 0.5806295405355471
$ws.Cells.Item(21, 7).Value = 0.284200842643159
$ws.Cells.Item(21, 8).Value = 0.4278422169529676
$ws.Cells.Item(21, 10).Value = 0.3136652770592434
$ws.Cells.Item(21, 13).Value = 8.866237085205
$ws.Cells.Item(21, 15).Value = 1.347551108144756

$ws.Cells.Item(22, 2).Value = 0.1997751957783578
$ws.Cells.Item(22, 4).Value = 0.1771421323318592
$ws.Cells.Item(22, 5).Value = 0.1956733521465139
$ws.Cells.Item(22, 6).Value = 0.5698194416559303
$ws.Cells.Item(22, 7).Value = 0.2845147090363014
$ws.Cells.Item(22, 8).Value = 0.4234630553584253
$ws.Cells.Item(22, 10).Value = 0.326984776037051
$ws.Cells.Item(22, 13).Value = 9.482817128851991
$ws.Cells.Item(22, 15).Value = 1.339321472460881

$ws.Cells.Item(23, 2).Value = 0.1948655829996682
$ws.Cells.Item(23, 4).Value = 0.1788894488175572
$ws.Cells.Item(23, 5).Value = 0.1938805885477279
$ws.Cells.Item(23, 6).Value = 0.5755052791164772
$ws.Cells.Item(23, 7).Value = 0.2842992702970122
$ws.Cells.Item(23, 8).Value = 0.4257634955027356
$ws.Cells.Item(23, 10).Value = 0.3198513048360212
$ws.Cells.Item(23, 13).Value = 9.153883446081863
$ws.Cells.Item(23, 15).Value = 1.343536074859117

$ws.Cells.Item(24, 2).Value = 0.1763896569628827
$ws.Cells.Item(24, 4).Value = 0.1858926162091237
$ws.Cells.Item(24, 5).Value = 0.1875204976867337
$ws.Cells.Item(24, 6).Value = 0.5988640312728108
$ws.Cells.Item(24, 7).Value = 0.2845255164344564
$ws.Cells.Item(24, 8).Value = 0.4352801552565069
$ws.Cells.Item(24, 10).Value = 0.2933674467256111
$ws.Cells.Item(24, 13).Value = 7.904853119674385
$ws.Cells.Item(24, 15).Value = 1.363370601002515

$ws.Cells.Item(25, 2).Value = 0.1567173484472164
$ws.Cells.Item(25, 4).Value = 0.1942574310576006
$ws.Cells.Item(25, 5).Value = 0.1815515617692753
$ws.Cells.Item(25, 6).Value = 0.62786579544386
$ws.Cells.Item(25, 7).Value = 0.2869288908642034
$ws.Cells.Item(25, 8).Value = 0.4472356450834525
$ws.Cells.Item(25, 10).Value = 0.2658934879451493
$ws.Cells.Item(25, 13).Value = 6.55052388243513
$ws.Cells.Item(25, 15).Value = 1.392837692152881
